$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '63.529.24'
$ws.Range('E2').Value = '  +1.58%  '

$ws.Range('D3').Value = '2.657.86'
$ws.Range('E3').Value = '  +3.52%  '

$ws.Range('D4').NumberFormat = '@'
$ws.Range('D4').Value = '0.998'
$ws.Range('E4').Value = '  -0.20%  '

$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '594.95'
$ws.Range('E5').Value = '  +2.97%  '

$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '144.28'
$ws.Range('E6').Value = '  +0.45%  '

$ws.Range('D7').NumberFormat = '@'
$ws.Range('D7').Value = '0.998'
$ws.Range('E7').Value = '  -0.22%  '

$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '0.588'
$ws.Range('E8').Value = '  +0.08%  '

$ws.Range('D9').Value = '2.659.98'
$ws.Range('E9').Value = '  +3.73%  '

$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '0.107'
$ws.Range('E10').Value = '  +0.87%  '

$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '5.69'
$ws.Range('E11').Value = '  +2.64%  '

$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '0.153'
$ws.Range('E12').Value = '  +0.91%  '

$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '0.356'
$ws.Range('E13').Value = '  +1.49%  '

$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '27.57'
$ws.Range('E14').Value = '  +2.74%  '

$ws.Range('D15').Value = '3.121.96'
$ws.Range('E15').Value = '  +3.11%  '

$ws.Range('D16').Value = '63.314.92'
$ws.Range('E16').Value = '  +1.34%  '

$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '0.0000145'
$ws.Range('E17').Value = '  +1.05%  '

$ws.Range('D18').Value = '2.652.05'
$ws.Range('E18').Value = '  +3.57%  '

$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '11.43'
$ws.Range('E19').Value = '  +2.51%  '

$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '340.89'
$ws.Range('E20').Value = '  +0.91%  '

$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '4.39'
$ws.Range('E21').Value = '  +1.42%  '

$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '6.78'
$ws.Range('E22').Value = '  +2.09%  '

$ws.Range('E23').Value = '  +0.04%  '

$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '67.13'
$ws.Range('E24').Value = '  +0.20%  '

$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '1.68'
$ws.Range('E25').Value = '  +6.96%  '

$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '1.55'
$ws.Range('E26').Value = '  +3.55%  '

$ws.Range('E27').Value = '  +1.03%  '

$ws.Range('B28').Value = 'InternetComputer(DFINITY)'
$ws.Range('C28').Value = 'https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp'
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '8.51'
$ws.Range('E28').Value = '  +3.72%  '

$ws.Range('B29').Value = 'Binance-PegBSC-USD'
$ws.Range('C29').Value = 'https://coinranking.com/coin/i5jggxiwp+binance-pegbsc-usd-bsc-usd'
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '1.00'
$ws.Range('E29').Value = '  -0.05%  '

$ws.Range('B30').Value = 'Bittensor'
$ws.Range('C30').Value = 'https://coinranking.com/coin/pgv7xSFi6+bittensor-tao'
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '538.76'
$ws.Range('E30').Value = '  +17.97%  '

$ws.Range('B31').Value = 'Aptos'
$ws.Range('C31').Value = 'https://coinranking.com/coin/HGYj5JCv5+aptos-apt'
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '7.86'
$ws.Range('E31').Value = '  -1.38%  '

$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '1.85'
$ws.Range('E32').Value = '  +15.05%  '

$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '1.99'
$ws.Range('E33').Value = '  +3.97%  '

$ws.Range('D34').Value = '0.0₃0815'
$ws.Range('E34').Value = '  +1.72%  '

$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '174.51'
$ws.Range('E35').Value = '  -0.95%  '

$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '5.04'
$ws.Range('E36').Value = '  +13.79%  '

$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '0.999'
$ws.Range('E37').Value = '  +0.06%  '

$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '0.404'
$ws.Range('E38').Value = '  +2.36%  '

$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '19.11'
$ws.Range('E39').Value = '  +1.41%  '

$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '1.82'
$ws.Range('E40').Value = '  +8.76%  '

$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '173.14'
$ws.Range('E41').Value = '  +9.00%  '

$ws.Range('E42').Value = '  -0.02%  '

$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '40.15'
$ws.Range('E43').Value = '  +0.43%  '

$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '3.77'
$ws.Range('E44').Value = '  +2.44%  '

$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '22.29'
$ws.Range('E45').Value = '  +7.13%  '

$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '0.0564'
$ws.Range('E46').Value = '  +5.85%  '

$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '0.633'
$ws.Range('E47').Value = '  +0.95%  '

$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '0.0963'
$ws.Range('E48').Value = '  +0.55%  '

$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '0.0240'
$ws.Range('E49').Value = '  +2.74%  '

$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '18.75'
$ws.Range('E50').Value = '  +4.48%  '

$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '1.72'
$ws.Range('E51').Value = '  +2.95%  '
